$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$DATE_VAL = 45796.38566594907

function Set-Row($r, $name, $email, $score) {
    $ws.Range("A$r").Value = $name
    $ws.Range("B$r").Value = $email
    $ws.Range("C$r").Value = $score
    $ws.Range("D$r").Value = $DATE_VAL
}

# Row 2-14 data (name, email, score)
Set-Row 2  "Test"              "b@bbbbbb.com"       300
Set-Row 3  "SSSSSSS"           "s@s.com"            280
Set-Row 4  "Test3"             "N/A"                20
Set-Row 5  "Test3"             "N/A"                20
Set-Row 6  "Test3"             "N/A"                20
Set-Row 7  "ADSGADFASDFASDFSA" "a@d.com"             20
Set-Row 8  "Test"              "N/A"                10
Set-Row 9  "Test"              "N/A"                10
Set-Row 10 "Test"              "N/A"                10
Set-Row 11 "Winner"            "aga@aga.com"         9
Set-Row 12 "sK"                "aaaa@a.xom"          0
Set-Row 13 "ghjgfjhfgjh"       "vvfhgfh@h.com"       0
Set-Row 14 "D"                 "d@d.com"             0

# Apply the date number format once, then clone that exact style onto the
# rest of column D via copy/paste-special so every row shares ONE style
# index (s=5) instead of the engine minting a brand new cellXfs entry for
# every NumberFormat assignment.
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Copy()
$ws.Range("D3:D14").PasteSpecial(-4122)

# Make sure blank E cells for the newly-added rows 8-14 carry the same
# style as the rest of column E (reuse existing font so no new style is
# created).
foreach ($r in 8..14) {
    $ws.Range("E$r").Font.Name = "Calibri"
}

# Row 11 / column B: "Winner"'s email is a manually-formatted hyperlink
# (underlined, blue) with a real mailto hyperlink attached.
$chars1 = $ws.Range("B11").Characters(1, 10)
$font1 = $chars1.Font
$font1.Underline = $true
$font1.ColorIndex = 5
$font1.Name = "Calibri"
$font1.Size = 11
$chars2 = $ws.Range("B11").Characters(11, 1)
$font2 = $chars2.Font
$font2.Underline = $true
$font2.ColorIndex = 5
$font2.Name = "Calibri"
$font2.Size = 11

$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:aga@aga.com")
# Hyperlinks.Add re-styles the whole cell with the built-in Hyperlink
# style; restore the plain per-row style (s=1) while keeping the rich
# text run formatting + the hyperlink relationship intact.
$ws.Range("A11").Copy()
$ws.Range("B11").PasteSpecial(-4122)

# Row heights
$ws.Rows.Item(2).RowHeight = 14.7
foreach ($r in 3..14) {
    $ws.Rows.Item($r).RowHeight = 13.55
}
